$wb = $excel.ActiveWorkbook

# 1) Rename the second worksheet's tab: "Include from RoleClass" -> "Include #0"
$wsInclude = $wb.Worksheets.Item("Include from RoleClass")
$wsInclude.Name = "Include #0"

$wsMeta = $wb.Worksheets.Item("Metadata")

# 2) Bump the Version value (row 3, column B)
$wsMeta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# 3) Bump the Date value (row 8, column B)
$wsMeta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# 4) Insert a new "Jurisdiction" property row right after "Contact" (row 10),
#    before "Description" (old row 11) -- this pushes Description/Purpose/
#    Copyright/Immutable down by one row.
$wsMeta.Rows.Item(11).Insert()

# Match the formatting of the surrounding property rows (border + wrap style)
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)

$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""

Write-Output "done"
